# Alter raw data to include categorical variable: replace the
# "First Letter of Cat's Name" column in the Data sheet with a
# "Cat Breed" column, and update the corresponding row in the
# Codebook sheet to describe the new variable.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$codebookSheet = $wb.Worksheets.Item("Codebook")

# --- Data sheet: column D ---
$dataSheet.Range("D1").Value = "Cat Breed"

$breeds = @("S", "S", "AS", "P", "P", "S", "AS", "MC", "M", "AS", "AS", "AS", "AS", "MC")
for ($i = 0; $i -lt $breeds.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 4).Value = $breeds[$i]
}

# --- Codebook sheet: row 5 describes the new "Cat Breed" variable ---
$codebookSheet.Range("A5").Value = "Cat Breed"
$codebookSheet.Range("B5").Value = "Breed of the individual's cat (Siamese /Maine Coon/American Shorthair/Persian)"
$codebookSheet.Range("C5").Value = "S/MC/AS/P/NA"

# --- restore selections on both sheets (Codebook stays the active tab) ---
$dataSheet.Activate() | Out-Null
$dataSheet.Range("G20").Select() | Out-Null

$codebookSheet.Activate() | Out-Null
$codebookSheet.Range("F8").Select() | Out-Null
